# Adds 18 new proverb rows (92-109) to Sheet1, sourced from "لال ءُ یاقوت",
# then pads four trailing blank (but centre-styled) rows (110-113), matching
# the row-by-row insertion order the workbook was originally authored in so
# that the shared-string table comes out in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entries = @(
    @{ Addr = "B92"; Text = "بےعزّت ہر جاہ بے عزّت اِنت آنہ وتی عزّت ءَ زانت نیکہ دگرءِ" },
    @{ Addr = "A92"; Text = "بے عزّت دگر ءَ بے عزّت کنت۔" },
    @{ Addr = "C92"; Text = "لال ءُ یاقوت" },
    @{ Addr = "B93"; Text = "کسے کہ کار نہ کنت، جہد نہ کنت بے بہر اِنت، پَشت کپیت" },
    @{ Addr = "A93"; Text = "بےکاری، بے بہری۔" },
    @{ Addr = "C93"; Text = "لال ءُ یاقوت" },
    @{ Addr = "B94"; Text = "بےکاریں مردم ہچیز ءِتہانپ نہ گندیت" },
    @{ Addr = "A94"; Text = "بےکاری، نَپَ نہ داری۔" },
    @{ Addr = "C94"; Text = "لال ءُ یاقوت" },
    @{ Addr = "B95"; Text = "بےکمال ءَ را ہچیز ءِ پرواہ نہ بیت" },
    @{ Addr = "A95"; Text = "بےکمال نہ سیال گندیت نہ مال۔" },
    @{ Addr = "C95"; Text = "لال ءُ یاقوت" },
    @{ Addr = "B96"; Text = "اگاں کسے ترا لوٹیت دعوت کنت گڑا برو، بے لوٹگ ءَ وتارابے شرپ مہ کن" },
    @{ Addr = "A96"; Text = "بےلوٹگ ءَ کسی نان ءَ مرو۔" },
    @{ Addr = "C96"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A97"; Text = "بےمال مال نہ کٹیت۔" },
    @{ Addr = "B97"; Text = "بے جوہر ءُ تاوان ہچ کرت نہ کنت گڑا مال ءَ کجا چہ کٹیت" },
    @{ Addr = "C97"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A98"; Text = "بے مال ءَ را سیال نہ بیت۔" },
    @{ Addr = "B98"; Text = "بے مال ءُ گریب ءَ کس وتی سیال نہ کنت، دنیا لالچی اِنت" },
    @{ Addr = "C98"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A99"; Text = "بے مالی بدحالی۔" },
    @{ Addr = "B99"; Text = "آئی ءَ کہ مال نہ بوت بزاں ہچ نہ بوت، بدحالی اِنت" },
    @{ Addr = "C99"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A100"; Text = "بے مالی بدحالی۔" },
    @{ Addr = "B100"; Text = "مال کہ نہ بوت بدحال بئے" },
    @{ Addr = "C100"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A101"; Text = "بے میاری پہ بلوچ ءَ نبرازی۔" },
    @{ Addr = "B101"; Text = "بلوچ میار جلیں ننگ داریں راجے آئی واستہ بے میاری عیب اِنت" },
    @{ Addr = "C101"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A102"; Text = "بے میار ءَ راپَیزنبی۔" },
    @{ Addr = "B102"; Text = "بے لج ءُ ننگ ءِ درور آپ ءِ لکیر کشگ" },
    @{ Addr = "C102"; Text = "لال ءُ یاقوت" },
    @{ Addr = "B103"; Text = "بے ایمان بے نان بیت" },
    @{ Addr = "A103"; Text = "بےنان بےایمان بیت۔" },
    @{ Addr = "C103"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A104"; Text = "بے نمازی بے روزی، تنگ نمازی تنگ روزی۔" },
    @{ Addr = "B104"; Text = "کسے ءِ تہا کہ دین نہ بیت، نماز نہ کنت گڑا آئی کِرّا‌ روزی نئیت تنگ بیت" },
    @{ Addr = "C104"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A105"; Text = "بے واہگ ءِ ہمسائگ مہ بو۔" },
    @{ Addr = "B105"; Text = "اَچ کسے ءَ چہ تماہ نہ گند ئے، واہشت نہ گند ئے وتارا بے شرف مہ کن پہ آئی مرہ، بے واہگیں مردم ءِ همسایگی اوں وش نہ اِنت" },
    @{ Addr = "C105"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A106"; Text = "بے واہگ پہ گلِگ ءَ جوان نہ بیت۔" },
    @{ Addr = "B106"; Text = "آکہ بے واہگ ءُ ناامیت اِنت آچ آئی ءَ گلِگ کنگ بے نپ اِنت" },
    @{ Addr = "C106"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A107"; Text = "بے وفا نبی سپا۔" },
    @{ Addr = "B107"; Text = "غدار ءُ دھوکہ باز، ہچبر نہ بنت پارسا" },
    @{ Addr = "C107"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A108"; Text = "بے وفا نگندی نپ ءَ۔" },
    @{ Addr = "B108"; Text = "آمردم کہ دغاباز اِنت بے وفا اِنت آ ہچبر سیت ءُ نپ نہ گندیت" },
    @{ Addr = "C108"; Text = "لال ءُ یاقوت" },
    @{ Addr = "A109"; Text = "بے ہمتی، بد قسمتی۔" },
    @{ Addr = "B109"; Text = "کسے کہ جہد نہ کنت، ہمت نہ کنت گڑا شومی کجام بہ بیت" },
    @{ Addr = "C109"; Text = "لال ءُ یاقوت" }
)

foreach ($entry in $entries) {
    $ws.Range($entry.Addr).Value = $entry.Text
}

# New rows inherit the same centred style ("s=2") used by every other data
# row on the sheet, including the still-empty buffer rows beneath the data.
$ws.Range("A92:C113").HorizontalAlignment = -4108

# Leave the freshly-entered last row selected, matching the author's cursor.
$ws.Range("C109").Select()

